$wb = $excel.ActiveWorkbook

$wsSize = $wb.Worksheets.Item("Size")
$wsCost = $wb.Worksheets.Item("Cost")
$wsInd = $wb.Worksheets.Item("Indicators")

$wsSize.Range("G2").Value = 320.37394354013548536
$wsSize.Range("G3").Value = 1620.94602417017199514
$wsSize.Range("G4").Value = 18.71482279474217947
$wsSize.Range("C5").Value = 244.98844740528031139
$wsSize.Range("D5").Value = 1069.25275700247698296
$wsSize.Range("E5").Value = 54.72222299471057738
$wsSize.Range("F5").Value = 47.44230428647126274
$wsSize.Range("G5").Value = 1416.40573168893911316

$wsCost.Range("I2").Value = 9.9840096059080814
$wsCost.Range("I3").Value = 0.64074788708027108
$wsCost.Range("I4").Value = 0.89152031329359449
$wsCost.Range("I5").Value = 0.003742964558948435
$wsCost.Range("E6").Value = 0.02449884474052803
$wsCost.Range("F6").Value = 0.10692527570024769
$wsCost.Range("G6").Value = 0.005472222299471058
$wsCost.Range("H6").Value = 0.004744230428647125
$wsCost.Range("I6").Value = 0.14164057316889389
$wsCost.Range("I7").Value = 0.1788789342984472
$wsCost.Range("I8").Value = 0.24888759957377349
$wsCost.Range("I9").Value = 0.005224656412621739
$wsCost.Range("E10").Value = 0.005129545481645616
$wsCost.Range("F10").Value = 0.02238791545686965
$wsCost.Range("G10").Value = 0.001145768850250167
$wsCost.Range("H10").Value = 0.000993342584799286
$wsCost.Range("I10").Value = 0.02965657237356472
$wsCost.Range("I11").Value = 0.13939342809793881
$wsCost.Range("E12").Value = 1.93345353896042194
$wsCost.Range("F12").Value = 5.66754634081156627
$wsCost.Range("G12").Value = 0.08265990018923595
$wsCost.Range("H12").Value = 0.0206568970888025
$wsCost.Range("I12").Value = 7.70431667705002621

$wsInd.Range("C2").Value = 57627.98830116898898268
$wsInd.Range("D2").Value = 167221.51660652519785799
$wsInd.Range("E2").Value = 224849.50490769420866854
$wsInd.Range("E3").Value = 0.76992550822477179
$wsInd.Range("E4").Value = 0.23007449177522829
$wsInd.Range("C5").Value = 0.82901061634023565
$wsInd.Range("D5").Value = 0.8787852344790047
$wsInd.Range("E5").Value = 0.86602820789371748
$wsInd.Range("E6").Value = 0.0002203891686185875
